$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "65.399.73"
$ws.Range("E2").Value = "  +1.62%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.649.11"
$ws.Range("E3").Value = "  +0.84%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.06%  "

# Row 5 - BNB
$ws.Range("D5").Value = "605.33"
$ws.Range("E5").Value = "  +1.54%  "

# Row 6 - Solana
# (156.10 would be coerced to the number 156.1 and lose its trailing zero,
# so force the cell to text first, then restore the default "Normal" style)
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.13%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.04%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -0.45%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.646.79"
$ws.Range("E9").Value = "  +0.82%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +7.46%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +1.36%  "

# Row 12 - Toncoin
$ws.Range("E12").Value = "  +0.10%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "29.85"
$ws.Range("E14").Value = "  +6.20%  "

# Row 15 - ShibaInu
$ws.Range("D15").Value = "0.0000195"
$ws.Range("E15").Value = "  +13.73%  "

# Row 16 - WrappedliquidstakedEther2.0
$ws.Range("D16").Value = "3.126.18"
$ws.Range("E16").Value = "  +0.96%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "65.191.61"
$ws.Range("E17").Value = "  +1.42%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.654.65"
$ws.Range("E18").Value = "  +1.14%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "12.72"
$ws.Range("E19").Value = "  +3.49%  "

# Row 20 - Polkadot
$ws.Range("E20").Value = "  +2.25%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "358.66"
$ws.Range("E21").Value = "  +2.69%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "7.47"
$ws.Range("E22").Value = "  +5.13%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.18%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "69.73"
$ws.Range("E24").Value = "  +3.00%  "

# Row 25 - SuiNetwork
$ws.Range("E25").Value = "  -0.20%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("D26").Value = "9.42"
$ws.Range("E26").Value = "  +1.66%  "

# Row 27 - PEPE
$ws.Range("D27").Value = "0.0000105"
$ws.Range("E27").Value = "  +14.82%  "

# Row 28 - Fetch.AI
$ws.Range("E28").Value = "  -2.62%  "

# Row 29 - Kaspa
$ws.Range("D29").Value = "0.166"
$ws.Range("E29").Value = "  +2.15%  "

# Row 30 - Aptos
$ws.Range("D30").Value = "8.11"
$ws.Range("E30").Value = "  -2.82%  "

# Row 31 - Binance-PegBSC-USD
# (1.00 would be coerced to the number 1, so force text as above)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.24%  "

# Row 32 - PancakeSwap
$ws.Range("D32").Value = "2.17"
$ws.Range("E32").Value = "  +4.17%  "

# Row 33 - Bittensor
$ws.Range("D33").Value = "522.87"
$ws.Range("E33").Value = "  -4.60%  "

# Row 34 - ImmutableX
$ws.Range("D34").Value = "1.78"
$ws.Range("E34").Value = "  -2.82%  "

# Row 35 - NEARProtocol
$ws.Range("E35").Value = "  -0.23%  "

# Row 36 - RenderToken
$ws.Range("E36").Value = "  +1.51%  "

# Row 37 - PolygonEcosystemToken
$ws.Range("E37").Value = "  +2.30%  "

# Row 38 - EthereumClassic
$ws.Range("D38").Value = "20.65"
$ws.Range("E38").Value = "  +2.73%  "

# Row 39 - Monero
$ws.Range("D39").Value = "162.31"
$ws.Range("E39").Value = "  -1.94%  "

# Row 40 - FirstDigitalUSD
$ws.Range("E40").Value = "  +0.05%  "

# Row 43 - OKB
$ws.Range("D43").Value = "41.95"
$ws.Range("E43").Value = "  +0.20%  "

# Row 44 - Aave
$ws.Range("D44").Value = "165.47"
$ws.Range("E44").Value = "  -1.83%  "

# Row 45 - Filecoin
$ws.Range("D45").Value = "4.12"
$ws.Range("E45").Value = "  -0.23%  "

# Row 46 - dogwifhat
$ws.Range("E46").Value = "  +4.08%  "

# Row 47 - Hedera
$ws.Range("D47").Value = "0.0608"
$ws.Range("E47").Value = "  +2.89%  "

# Row 48 - InjectiveProtocol
# (22.90 would be coerced to the number 22.9, so force text as above)
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.39%  "

# Row 49 - was Mantle, now VeChain
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "0.0263"
$ws.Range("E49").Value = "  +3.92%  "

# Row 50 - was VeChain, now Mantle
# (0.650 would be coerced to the number 0.65, so force text as above)
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.650"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.38%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  +0.22%  "
